$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: updated frequency values for publication
$ws.Range("B2").Value = 0.145515525262361
$ws.Range("C2").Value = 0.172130260737856
$ws.Range("D2").Value = 0.761224710591799
$ws.Range("E2").Value = 0.808070972627935
$ws.Range("F2").Value = 0.178910166251938
$ws.Range("G2").Value = 0.695589455083126
$ws.Range("H2").Value = 0.116520610191496
$ws.Range("I2").Value = 0.122939882433553
$ws.Range("J2").Value = 0.922824479786505
$ws.Range("K2").Value = 0.298568285910058
$ws.Range("L2").Value = 0.0973709834469328
$ws.Range("M2").Value = 0.790003245699448
$ws.Range("N2").Value = 0.219084712755599
$ws.Range("O2").Value = 0.967218435572866
$ws.Range("P2").Value = 0.374445526344261
$ws.Range("Q2").Value = 0.0968300335388943
$ws.Range("R2").Value = 0.966533232356017
$ws.Range("S2").Value = 0.823974899924267
$ws.Range("T2").Value = 0.784593746619063
$ws.Range("U2").Value = 0.780013704064337
$ws.Range("V2").Value = 0.0872371885030113
$ws.Range("W2").Value = 0.0995708464062895
$ws.Range("X2").Value = 0.131162321035739

# Row 3: updated frequency values for publication
$ws.Range("B3").Value = 0.065491002199863
$ws.Range("C3").Value = 0.106927765155613
$ws.Range("D3").Value = 0.0338995275704136
$ws.Range("E3").Value = 0.068989144938512
$ws.Range("F3").Value = 0.618413934869631
$ws.Range("G3").Value = 0.161996465793934
$ws.Range("H3").Value = 0.0387680767427603
$ws.Range("I3").Value = 0.0561866637816005
$ws.Range("J3").Value = 0.0209527931046918
$ws.Range("K3").Value = 0.0491903783043024
$ws.Range("L3").Value = 0.0164088138771683
$ws.Range("M3").Value = 0.0355223772945292
$ws.Range("N3").Value = 0.0680875617584478
$ws.Range("O3").Value = 0.00645533556925962
$ws.Range("P3").Value = 0.00717660211331097
$ws.Range("Q3").Value = 0.858595694038732
$ws.Range("R3").Value = 0.00894370514623679
$ws.Range("S3").Value = 0.0772476468679
$ws.Range("T3").Value = 0.0232969093728587
$ws.Range("U3").Value = 0.0355584406217318
$ws.Range("V3").Value = 0.07191027444192
$ws.Range("W3").Value = 0.0257492156226333
$ws.Range("X3").Value = 0.0850012622164521

# Row 4: updated frequency values for publication
$ws.Range("B4").Value = 0.67921670453316
$ws.Range("C4").Value = 0.0447906523855891
$ws.Range("D4").Value = 0.0693137148833351
$ws.Range("E4").Value = 0.0271556853835335
$ws.Range("F4").Value = 0.0424825994446248
$ws.Range("G4").Value = 0.0478199718706048
$ws.Range("H4").Value = 0.0842078690179956
$ws.Range("I4").Value = 0.0668974719607631
$ws.Range("J4").Value = 0.0256770889682282
$ws.Range("K4").Value = 0.588084676692272
$ws.Range("L4").Value = 0.0974070467741354
$ws.Range("M4").Value = 0.158065563128854
$ws.Range("N4").Value = 0.665837210141008
$ws.Range("O4").Value = 0.00768148869414692
$ws.Range("P4").Value = 0.605936023657543
$ws.Range("Q4").Value = 0.0208085397958816
$ws.Range("R4").Value = 0.00151465974250784
$ws.Range("S4").Value = 0.0921057376753579
$ws.Range("T4").Value = 0.0226477694832125
$ws.Range("U4").Value = 0.0223231995383894
$ws.Range("V4").Value = 0.0836308557827545
$ws.Range("W4").Value = 0.706913339824732
$ws.Range("X4").Value = 0.707129719787948

# Row 5: updated frequency values for publication
$ws.Range("B5").Value = 0.109704641350211
$ws.Range("C5").Value = 0.676043131739334
$ws.Range("D5").Value = 0.135453856972844
$ws.Range("E5").Value = 0.0957481337228173
$ws.Range("F5").Value = 0.160193299433806
$ws.Range("G5").Value = 0.0945580439251325
$ws.Range("H5").Value = 0.760503444047748
$ws.Range("I5").Value = 0.753831728515273
$ws.Range("J5").Value = 0.0305456381405749
$ws.Range("K5").Value = 0.064156659093368
$ws.Range("L5").Value = 0.788813155901763
$ws.Range("M5").Value = 0.0162284972411555
$ws.Range("N5").Value = 0.0469183886905406
$ws.Range("O5").Value = 0.0186447401637275
$ws.Range("P5").Value = 0.0123697212304807
$ws.Range("Q5").Value = 0.0237296692992896
$ws.Range("R5").Value = 0.0230084027552382
$ws.Range("S5").Value = 0.00656352555086732
$ws.Range("T5").Value = 0.169461574524866
$ws.Range("U5").Value = 0.162032529121137
$ws.Range("V5").Value = 0.757221681272314
$ws.Range("W5").Value = 0.167730534819142
$ws.Range("X5").Value = 0.076670633632659

